# Update "想去人数" (F column) counts across sheets, as produced by the
# latest data-refresh run (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1228
$ws1.Range("F3").Value = 433
$ws1.Range("F4").Value = 21
$ws1.Range("F5").Value = 12508
$ws1.Range("F6").Value = 68
$ws1.Range("F7").Value = 27
$ws1.Range("F9").Value = 6
$ws1.Range("F10").Value = 12388
$ws1.Range("F11").Value = 235
$ws1.Range("F12").Value = 4891
$ws1.Range("F13").Value = 4809
$ws1.Range("F15").Value = 72
$ws1.Range("F17").Value = 104
$ws1.Range("F18").Value = 961
$ws1.Range("F21").Value = 174
$ws1.Range("F23").Value = 5217

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 8

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1228
$ws4.Range("F3").Value = 433
$ws4.Range("F4").Value = 8
$ws4.Range("F6").Value = 21
$ws4.Range("F7").Value = 12508
$ws4.Range("F8").Value = 68
$ws4.Range("F9").Value = 27
$ws4.Range("F11").Value = 6
$ws4.Range("F12").Value = 12388
$ws4.Range("F13").Value = 235
$ws4.Range("F14").Value = 4891
$ws4.Range("F15").Value = 4809
$ws4.Range("F17").Value = 72
$ws4.Range("F19").Value = 104
$ws4.Range("F20").Value = 961
$ws4.Range("F23").Value = 174
$ws4.Range("F25").Value = 5217
